$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18, shifting existing rows 18:112 down to 19:113
$ws.Rows.Item(18).Insert()

# Populate the new row 18 with data (copy of unchanged columns + new values)
$ws.Cells.Item(18, 1).Value = 3
$ws.Cells.Item(18, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(18, 3).Value = "Coquimbo"
$ws.Cells.Item(18, 4).NumberFormat = $ws.Cells.Item(19, 4).NumberFormat
$ws.Cells.Item(18, 4).Value = 44550
$ws.Cells.Item(18, 5).Value = 5
$ws.Cells.Item(18, 6).Value = 100112030
$ws.Cells.Item(18, 7).Value = "Poroto granado"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 35
$ws.Cells.Item(18, 11).Value = 45000
$ws.Cells.Item(18, 12).Value = 45000
$ws.Cells.Item(18, 13).Value = 45000
$ws.Cells.Item(18, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(18, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(18, 16).Value = 1800
$ws.Cells.Item(18, 17).Value = 25
$ws.Cells.Item(18, 18).Value = "Hortaliza"
